$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.282.14'
$ws.Range('D3').Value = '2.250.66'
$ws.Range('E3').Value = '  +1.10%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '''307.66'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.57%  '
$ws.Range('D6').Value = '''96.62'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.27%  '
$ws.Range('D7').Value = '''0.574'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +1.30%  '
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('D9').Value = '''0.529'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.60%  '
$ws.Range('D10').Value = '''35.35'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.07%  '
$ws.Range('D11').Value = '''0.0817'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('D12').Value = '''7.29'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.77%  '
$ws.Range('D13').Value = '''0.105'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range('D14').Value = '2.589.41'
$ws.Range('E14').Value = '  +1.03%  '
$ws.Range('D15').Value = '2.323.90'
$ws.Range('E15').Value = '  +4.68%  '
$ws.Range('D16').Value = '''0.839'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.30%  '
$ws.Range('D17').Value = '''13.66'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.48%  '
$ws.Range('D18').Value = '44.116.86'
$ws.Range('E18').Value = '  +1.35%  '
$ws.Range('D19').Value = '0.0₃0973'
$ws.Range('E19').Value = '  +1.48%  '
$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D20').Value = '''12.26'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -6.18%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '''6.41'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.98%  '
$ws.Range('D22').Value = '''65.70'
$ws.Range('D22').ClearFormats()
$ws.Range('D23').Value = '''237.61'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.06%  '
$ws.Range('D24').Value = '''2.98'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.73%  '
$ws.Range('D25').Value = '''2.02'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.22%  '
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('D27').Value = '''10.05'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.78%  '
$ws.Range('D28').Value = '''2.21'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.59%  '
$ws.Range('D29').Value = '''37.98'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +3.95%  '
$ws.Range('D30').Value = '''5.99'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.41%  '
$ws.Range('D31').Value = '''20.22'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.10%  '
$ws.Range('D32').Value = '''153.11'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.88%  '
$ws.Range('D33').Value = '''0.0804'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.59%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = '''3.24'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +4.15%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').Value = '''2.61'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.86%  '
$ws.Range('D36').Value = '''0.121'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +3.22%  '
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('D38').Value = '''1.77'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -5.21%  '
$ws.Range('D39').Value = '''3.58'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.52%  '
$ws.Range('D40').Value = '''14.69'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -5.90%  '
$ws.Range('D41').Value = '''3.88'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.66%  '
$ws.Range('D42').Value = '''0.0300'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.93%  '
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('D44').Value = '1.760.20'
$ws.Range('E44').Value = '  +3.74%  '
$ws.Range('D45').Value = '''83.51'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('D46').Value = '''0.193'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.19%  '
$ws.Range('D47').Value = '''100.89'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.47%  '
$ws.Range('D48').Value = '''4.97'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.58%  '
$ws.Range('D49').Value = '''8.21'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.13%  '
$ws.Range('D50').Value = '''55.16'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.20%  '
$ws.Range('E51').Value = '  -3.12%  '
